$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update F column timestamps on the "data" sheet (rows 2-72) ---
$newTimestamps = @(
    "2021-10-05 14:34:08.326918",
    "2021-10-05 14:34:08.326926",
    "2021-10-05 14:34:08.326929",
    "2021-10-05 14:34:08.326932",
    "2021-10-05 14:34:08.326934",
    "2021-10-05 14:34:08.326937",
    "2021-10-05 14:34:08.326939",
    "2021-10-05 14:34:08.326942",
    "2021-10-05 14:34:08.326944",
    "2021-10-05 14:34:08.326947",
    "2021-10-05 14:34:08.326949",
    "2021-10-05 14:34:08.326951",
    "2021-10-05 14:34:08.326954",
    "2021-10-05 14:34:08.326956",
    "2021-10-05 14:34:08.326958",
    "2021-10-05 14:34:08.326961",
    "2021-10-05 14:34:08.326963",
    "2021-10-05 14:34:08.326966",
    "2021-10-05 14:34:08.326968",
    "2021-10-05 14:34:08.326971",
    "2021-10-05 14:34:08.326973",
    "2021-10-05 14:34:08.326976",
    "2021-10-05 14:34:08.326978",
    "2021-10-05 14:34:08.326980",
    "2021-10-05 14:34:08.326983",
    "2021-10-05 14:34:08.326985",
    "2021-10-05 14:34:08.326988",
    "2021-10-05 14:34:08.326990",
    "2021-10-05 14:34:08.326993",
    "2021-10-05 14:34:08.326995",
    "2021-10-05 14:34:08.326999",
    "2021-10-05 14:34:08.327002",
    "2021-10-05 14:34:08.327004",
    "2021-10-05 14:34:08.327007",
    "2021-10-05 14:34:08.327009",
    "2021-10-05 14:34:08.327012",
    "2021-10-05 14:34:08.327014",
    "2021-10-05 14:34:08.327016",
    "2021-10-05 14:34:08.327019",
    "2021-10-05 14:34:08.327021",
    "2021-10-05 14:34:08.327024",
    "2021-10-05 14:34:08.327026",
    "2021-10-05 14:34:08.327029",
    "2021-10-05 14:34:08.327031",
    "2021-10-05 14:34:08.327034",
    "2021-10-05 14:34:08.327036",
    "2021-10-05 14:34:08.327038",
    "2021-10-05 14:34:08.327041",
    "2021-10-05 14:34:08.327043",
    "2021-10-05 14:34:08.327045",
    "2021-10-05 14:34:08.327048",
    "2021-10-05 14:34:08.327050",
    "2021-10-05 14:34:08.327053",
    "2021-10-05 14:34:08.327055",
    "2021-10-05 14:34:08.327058",
    "2021-10-05 14:34:08.327060",
    "2021-10-05 14:34:08.327062",
    "2021-10-05 14:34:08.327065",
    "2021-10-05 14:34:08.327067",
    "2021-10-05 14:34:08.327069",
    "2021-10-05 14:34:08.327072",
    "2021-10-05 14:34:08.327074",
    "2021-10-05 14:34:08.327077",
    "2021-10-05 14:34:08.327079",
    "2021-10-05 14:34:08.327082",
    "2021-10-05 14:34:08.327085",
    "2021-10-05 14:34:08.327087",
    "2021-10-05 14:34:08.327090",
    "2021-10-05 14:34:08.327092",
    "2021-10-05 14:34:08.327094",
    "2021-10-05 14:34:08.327097"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# --- Add the "metadata" worksheet right after "data" ---
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$ws.Name = "metadata"

# Header row (B1:G1) -- bold, bordered, centered, matches the "data" sheet's header style
$headerRange = $ws.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Index cell (A2) -- same bold/border/centered style as the header row
$aCell = $ws.Range("A2")
$aCell.Value = 0
$aCell.Font.Bold = $true
$aCell.Borders.LineStyle = 1
$aCell.HorizontalAlignment = -4108
$aCell.VerticalAlignment = -4160

$ws.Range("B2").Value = "Immune_markers_WTS_UMCCR"
$ws.Range("C2").Value = 243

# data_version must stay textual ("0.75"), not be coerced to a number
$dVersion = $ws.Range("D2")
$dVersion.NumberFormat = "@"
$dVersion.Value = "0.75"
$dVersion.ClearFormats()

$ws.Range("E2").Value = "2019-11-22T04:11:29.121287Z"
$ws.Range("F2").Value = "2021-10-05 14:34:08.323769"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/243/?format=json"

# Keep "data" as the active sheet/tab (unchanged in the target workbook view)
$dataSheet.Activate() | Out-Null
$dataSheet.Range("A1").Select() | Out-Null
